$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest quarter column (D): this shifts every column one slot to the
# left (old E -> D, old F -> E, ... old M -> L) and frees up column M for the
# newest quarter being added.
$ws.Range("D:D").Delete() | Out-Null

# --- Column M: the newly published quarter -------------------------------

# Header row (quarter label) and the "as of" date row. Plain text that looks
# like a bare Excel date ("1402-02-27") gets auto-coerced to a date serial by
# .Value, so round-trip it through a throwaway formula cell + paste-values
# instead of typing it directly (keeps it a literal shared string, no
# unwanted numeric conversion and no NumberFormat-driven new style).
$ws.Range("O1").Formula = "=""1402-02-27"""
$ws.Range("O1").Copy() | Out-Null
$ws.Range("M9").PasteSpecial(-4163) | Out-Null
$ws.Range("O1").ClearContents() | Out-Null

$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Financial figures for the new quarter.
$ws.Range("M11").Value = 181187
$ws.Range("M12").Value = -186074
$ws.Range("M13").Value = -4887
$ws.Range("M14").Value = -3217
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 7590
$ws.Range("M17").Value = -514
$ws.Range("M18").Value = -911
$ws.Range("M19").Value = 12692
$ws.Range("M20").Value = 11267
$ws.Range("M21").Value = 4463
$ws.Range("M22").Value = 15730
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 15730
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 17601
$ws.Range("M27").Value = 0

# The column-D delete left column M without the formatting that used to live
# on column M (it shifted away with the data); re-apply it from column L so
# M keeps the same style ids as the rest of the table.
$ws.Range("L1:L28").Copy() | Out-Null
$ws.Range("M1:M28").PasteSpecial(-4122) | Out-Null

# --- Republish one of the surviving date labels ---------------------------
# Column I (old column J) shows a revised publish date/revision marker.
$ws.Range("I9").Value = "1402-02-27 (7)"

$excel.CutCopyMode = 0
